$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# hunk 0: row 5
$ws.Range("H5").Value = 186.71428
$ws.Range("I5").Value = 186.71428
$ws.Range("K5").Value = 186.71428
$ws.Range("M5").Value = -71.71428
# hunk 1: row 28
$ws.Range("H28").Value = 1723.2222
$ws.Range("I28").Value = 1661
$ws.Range("K28").Value = 1661
$ws.Range("M28").Value = -1176
# hunk 2: row 80
$ws.Range("H80").Value = 100010900
$ws.Range("I80").Value = 142857800
$ws.Range("J80").Value = 34808.668
$ws.Range("K80").Value = 428573400
$ws.Range("L80").Value = 104426.004
$ws.Range("M80").Value = -428572402
$ws.Range("N80").Value = -106422.004
# hunk 3: row 83
$ws.Range("H83").Value = 100010900
$ws.Range("I83").Value = 142857800
$ws.Range("J83").Value = 34808.668
$ws.Range("K83").Value = 1285720200
$ws.Range("L83").Value = 313278.012
$ws.Range("M83").Value = -1285715208
$ws.Range("N83").Value = -323262.012
# hunk 4: row 86
$ws.Range("H86").Value = 222223340
$ws.Range("J86").Value = 500000740
$ws.Range("L86").Value = 500000740
$ws.Range("N86").Value = -500002986
# hunk 5: row 89
$ws.Range("H89").Value = 222223340
$ws.Range("J89").Value = 500000740
$ws.Range("L89").Value = 2500003700
$ws.Range("N89").Value = -2500014932
# hunk 6: row 96
$ws.Range("H96").Value = 1340.1428
$ws.Range("J96").Value = 1651
$ws.Range("L96").Value = 4953
$ws.Range("N96").Value = -7699
# hunk 7: row 137
$ws.Range("H137").Value = 2821422.5
$ws.Range("J137").Value = 6947829
$ws.Range("L137").Value = 20843487
$ws.Range("N137").Value = -20848587
# hunk 8: row 141
$ws.Range("H141").Value = 4947.2705
$ws.Range("I141").Value = 3911.7856
$ws.Range("J141").Value = 8168.778
$ws.Range("K141").Value = 11735.3568
$ws.Range("L141").Value = 24506.334
$ws.Range("M141").Value = -6555.356800000001
$ws.Range("N141").Value = -34866.334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# hunk 9: row 32
$ws.Range("H32").Value = 17544938
$ws.Range("I32").Value = 18519322
$ws.Range("K32").Value = 18519322
$ws.Range("M32").Value = -18519035

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# hunk 10: row 86
$ws.Range("H86").Value = 1349.2
$ws.Range("I86").Value = 1640.8572
$ws.Range("J86").Value = 668.6667
$ws.Range("K86").Value = 1640.8572
$ws.Range("L86").Value = 668.6667
$ws.Range("M86").Value = -517.8571999999999
$ws.Range("N86").Value = -2914.6667
# hunk 11: row 89
$ws.Range("H89").Value = 1349.2
$ws.Range("I89").Value = 1640.8572
$ws.Range("J89").Value = 668.6667
$ws.Range("K89").Value = 8204.286
$ws.Range("L89").Value = 3343.3335
$ws.Range("M89").Value = -2588.286
$ws.Range("N89").Value = -14575.3335
# hunk 12: row 96
$ws.Range("H96").Value = 72666
$ws.Range("I96").Value = 52500
$ws.Range("J96").Value = 112998
$ws.Range("K96").Value = 52500
$ws.Range("L96").Value = 112998
$ws.Range("M96").Value = -49754
$ws.Range("N96").Value = -118490
# hunk 13: row 99
$ws.Range("H99").Value = 2132.5
$ws.Range("I99").Value = 1198.75
$ws.Range("K99").Value = 1198.75
$ws.Range("M99").Value = 299.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# hunk 14: row 31
$ws.Range("H31").Value = 5087.775
$ws.Range("I31").Value = 3444.375
$ws.Range("K31").Value = 3444.375
$ws.Range("M31").Value = -3149.375
# hunk 15: row 34
$ws.Range("H34").Value = 5087.775
$ws.Range("I34").Value = 3444.375
$ws.Range("K34").Value = 3444.375
$ws.Range("M34").Value = -3242.375
# hunk 16: row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# hunk 17: row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
# hunk 18: row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# hunk 19: row 38
$ws.Range("H38").Value = 45.61111
$ws.Range("I38").Value = 32.9375
$ws.Range("J38").Value = 147
$ws.Range("K38").Value = 98.8125
$ws.Range("L38").Value = 441
$ws.Range("M38").Value = 248.1875
$ws.Range("N38").Value = -1135
# hunk 20: row 92
$ws.Range("H92").Value = 1225.4
$ws.Range("J92").Value = 1281.8334
$ws.Range("L92").Value = 3845.5002
$ws.Range("N92").Value = -6341.5002
# hunk 21: row 100
$ws.Range("H100").Value = 69331.664
$ws.Range("J100").Value = 101997.5
$ws.Range("L100").Value = 305992.5
$ws.Range("N100").Value = -307614.5
# hunk 22: row 107
$ws.Range("H107").Value = 614.9167
$ws.Range("J107").Value = 899.75
$ws.Range("L107").Value = 2699.25
$ws.Range("N107").Value = -6539.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# hunk 23: row 93
$ws.Range("H93").Value = 29888
$ws.Range("J93").Value = 29888
$ws.Range("L93").Value = 29888
$ws.Range("N93").Value = -33632
# hunk 24: row 113
$ws.Range("H113").Value = 34048.1
$ws.Range("J113").Value = 36664
$ws.Range("L113").Value = 36664
$ws.Range("N113").Value = -41004
# hunk 25: row 117
$ws.Range("H117").Value = 75995
$ws.Range("J117").Value = 75995
$ws.Range("L117").Value = 75995
$ws.Range("N117").Value = -82879
# hunk 26: row 122
$ws.Range("H122").Value = 3937.6
$ws.Range("I122").Value = 3937.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11812.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9362.799999999999
$ws.Range("N122").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# hunk 27: row 22
$ws.Range("H22").Value = 1653.3334
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# hunk 28: row 27
$ws.Range("H27").Value = 1653.3334
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# hunk 29: row 46
$ws.Range("H46").Value = 3713.0715
$ws.Range("J46").Value = 3767.8462
$ws.Range("L46").Value = 3767.8462
$ws.Range("N46").Value = -4143.8462
# hunk 30: row 61
$ws.Range("H61").Value = 3333.3333
$ws.Range("I61").Value = 3250
$ws.Range("K61").Value = 3250
$ws.Range("M61").Value = -3048
# hunk 31: row 82
$ws.Range("H82").Value = 2101.3547
$ws.Range("J82").Value = 2244.2666
$ws.Range("L82").Value = 2244.2666
$ws.Range("N82").Value = -2966.2666
# hunk 32: row 85
$ws.Range("H85").Value = 2101.3547
$ws.Range("J85").Value = 2244.2666
$ws.Range("L85").Value = 2244.2666
$ws.Range("N85").Value = -4740.2666
# hunk 33: row 93
$ws.Range("H93").Value = 2613.077
$ws.Range("I93").Value = 2490
$ws.Range("K93").Value = 2490
$ws.Range("M93").Value = -1242
# hunk 34: row 113
$ws.Range("H113").Value = 3333.3333
$ws.Range("I113").Value = 3250
$ws.Range("K113").Value = 3250
$ws.Range("M113").Value = -1080

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# hunk 35: row 97
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
# hunk 36: row 116
$ws.Range("H116").Value = 119992.336
$ws.Range("J116").Value = 119992.336
$ws.Range("L116").Value = 119992.336
$ws.Range("N116").Value = -129170.336
# hunk 37: row 132
$ws.Range("H132").Value = 2497
$ws.Range("I132").Value = 2095.7273
$ws.Range("K132").Value = 6287.1819
$ws.Range("M132").Value = -3757.1819

Write-Output "Applied all Sheets edits"
